$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 2.257119139371683
$ws.Range("E4").Value = 0.8998700214674639
$ws.Range("C5").Value = 2.257119139371683
$ws.Range("D5").Value = 0.8998700214674639
